$wb = $excel.ActiveWorkbook

# --- 1. Merge "cliente" and "funcionario" data into "utilizador" ---
$util = $wb.Worksheets.Item("utilizador")
$cliente = $wb.Worksheets.Item("cliente")
$funcionario = $wb.Worksheets.Item("funcionario")

# Headers for new columns J, K, L (from cliente) and M, N (from funcionario)
$util.Range("J1").Value = "pontos"
$util.Range("K1").Value = "newsletter"
$util.Range("L1").Value = "pagamento"
$util.Range("M1").Value = "administrator"
$util.Range("N1").Value = "salario"

# Build lookup tables keyed by email
$clienteMap = @{}
for ($r = 2; $r -le 18; $r++) {
    $email = $cliente.Cells.Item($r, 1).Value2
    if ($email) {
        $clienteMap[$email] = @($cliente.Cells.Item($r, 2).Value2, $cliente.Cells.Item($r, 3).Value2, $cliente.Cells.Item($r, 4).Value2)
    }
}

$funcMap = @{}
for ($r = 2; $r -le 5; $r++) {
    $email = $funcionario.Cells.Item($r, 1).Value2
    if ($email) {
        $funcMap[$email] = @($funcionario.Cells.Item($r, 2).Value2, $funcionario.Cells.Item($r, 3).Value2)
    }
}

for ($r = 2; $r -le 22; $r++) {
    $email = $util.Cells.Item($r, 1).Value2
    if ($clienteMap.ContainsKey($email)) {
        $vals = $clienteMap[$email]
        $util.Cells.Item($r, 10).Value = $vals[0]
        $util.Cells.Item($r, 11).Value = $vals[1]
        $util.Cells.Item($r, 12).Value = $vals[2]
    }
    if ($funcMap.ContainsKey($email)) {
        $vals = $funcMap[$email]
        $util.Cells.Item($r, 13).Value = $vals[0]
        $util.Cells.Item($r, 14).Value = $vals[1]
    }
}

# --- 2. Delete the now-redundant "cliente" and "funcionario" sheets ---
# Delete in reverse tab order so removing one doesn't shift a still-pending
# worksheet reference onto a different sheet.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("funcionario").Delete()
$wb.Worksheets.Item("cliente").Delete()
$excel.DisplayAlerts = $true

# --- 3. Make "utilizador" the active sheet/tab ---
# Re-fetch by name: worksheet references captured before the deletes above
# can now resolve to a different sheet because tab positions shifted.
$util = $wb.Worksheets.Item("utilizador")
$util.Activate()
$util.Select()
